# Apply "Atualizacao de bases das ligas" update to Costa Rica Primera Division sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Part 1: swap odds/result data between row-pairs that share the same match date ---
# (ids in column A stay fixed; B and F:AC are swapped between the two rows)

# Rows 104 <-> 105
$ws.Range("B104").Value = 5965190
$ws.Range("F104").Value = "AD San Carlos"
$ws.Range("G104").Value = "Sporting San Jose"
$ws.Range("H104").Value = 1
$ws.Range("I104").Value = 1
$ws.Range("J104").Value = "D"
$ws.Range("K104").Value = 2.15
$ws.Range("L104").Value = 3.25
$ws.Range("M104").Value = 3.2
$ws.Range("N104").Value = 2.15
$ws.Range("O104").Value = 3.2
$ws.Range("P104").Value = 3.25
$ws.Range("Q104").Value = -0.25
$ws.Range("R104").Value = 1.925
$ws.Range("S104").Value = 1.875
$ws.Range("T104").Value = 2.5
$ws.Range("U104").Value = 1.875
$ws.Range("V104").Value = 1.925
$ws.Range("W104").Value = -1
$ws.Range("X104").Value = 2.2
$ws.Range("Y104").Value = -1
$ws.Range("Z104").Value = -0.5
$ws.Range("AA104").Value = 0.4375
$ws.Range("AB104").Value = -1
$ws.Range("AC104").Value = 0.925

$ws.Range("B105").Value = 5963936
$ws.Range("F105").Value = "Municipal Perez Zeledon"
$ws.Range("G105").Value = "Santos de Gupiles"
$ws.Range("H105").Value = 1
$ws.Range("I105").Value = 1
$ws.Range("J105").Value = "D"
$ws.Range("K105").Value = 1.95
$ws.Range("L105").Value = 3.25
$ws.Range("M105").Value = 3.75
$ws.Range("N105").Value = 2.15
$ws.Range("O105").Value = 3
$ws.Range("P105").Value = 3.5
$ws.Range("Q105").Value = -0.25
$ws.Range("R105").Value = 1.85
$ws.Range("S105").Value = 1.95
$ws.Range("T105").Value = 2
$ws.Range("U105").Value = 1.85
$ws.Range("V105").Value = 1.95
$ws.Range("W105").Value = -1
$ws.Range("X105").Value = 2
$ws.Range("Y105").Value = -1
$ws.Range("Z105").Value = -0.5
$ws.Range("AA105").Value = 0.475
$ws.Range("AB105").Value = 0
$ws.Range("AC105").Value = 0

# Rows 124 <-> 125
$ws.Range("B124").Value = 5965201
$ws.Range("F124").Value = "Herediano"
$ws.Range("G124").Value = "Municipal Perez Zeledon"
$ws.Range("H124").Value = 3
$ws.Range("I124").Value = 1
$ws.Range("J124").Value = "H"
$ws.Range("K124").Value = 1.333
$ws.Range("L124").Value = 4.75
$ws.Range("M124").Value = 7
$ws.Range("N124").Value = 1.285
$ws.Range("O124").Value = 5
$ws.Range("P124").Value = 8.5
$ws.Range("Q124").Value = -1.5
$ws.Range("R124").Value = 1.825
$ws.Range("S124").Value = 1.975
$ws.Range("T124").Value = 3
$ws.Range("U124").Value = 1.875
$ws.Range("V124").Value = 1.925
$ws.Range("W124").Value = 0.2849999999999999
$ws.Range("X124").Value = -1
$ws.Range("Y124").Value = -1
$ws.Range("Z124").Value = 0.825
$ws.Range("AA124").Value = -1
$ws.Range("AB124").Value = 0.875
$ws.Range("AC124").Value = -1

$ws.Range("B125").Value = 5963942
$ws.Range("F125").Value = "Sporting San Jose"
$ws.Range("G125").Value = "Deportivo Saprissa"
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = "D"
$ws.Range("K125").Value = 3.3
$ws.Range("L125").Value = 3.2
$ws.Range("M125").Value = 2
$ws.Range("N125").Value = 3.75
$ws.Range("O125").Value = 3.25
$ws.Range("P125").Value = 1.85
$ws.Range("Q125").Value = 0.5
$ws.Range("R125").Value = 1.9
$ws.Range("S125").Value = 1.9
$ws.Range("T125").Value = 2.5
$ws.Range("U125").Value = 1.85
$ws.Range("V125").Value = 1.95
$ws.Range("W125").Value = -1
$ws.Range("X125").Value = 2.25
$ws.Range("Y125").Value = -1
$ws.Range("Z125").Value = 0.8999999999999999
$ws.Range("AA125").Value = -1
$ws.Range("AB125").Value = -1
$ws.Range("AC125").Value = 0.95

# Rows 130 <-> 133
$ws.Range("B130").Value = 5965203
$ws.Range("F130").Value = "Cartagines"
$ws.Range("G130").Value = "Sporting San Jose"
$ws.Range("H130").Value = 3
$ws.Range("I130").Value = 2
$ws.Range("J130").Value = "H"
$ws.Range("K130").Value = 1.85
$ws.Range("L130").Value = 3.5
$ws.Range("M130").Value = 3.4
$ws.Range("N130").Value = 2.2
$ws.Range("O130").Value = 3.25
$ws.Range("P130").Value = 2.875
$ws.Range("Q130").Value = -0.25
$ws.Range("R130").Value = 1.975
$ws.Range("S130").Value = 1.825
$ws.Range("T130").Value = 2.5
$ws.Range("U130").Value = 1.925
$ws.Range("V130").Value = 1.875
$ws.Range("W130").Value = 1.2
$ws.Range("X130").Value = -1
$ws.Range("Y130").Value = -1
$ws.Range("Z130").Value = 0.9750000000000001
$ws.Range("AA130").Value = -1
$ws.Range("AB130").Value = 0.925
$ws.Range("AC130").Value = -1

$ws.Range("B133").Value = 5965205
$ws.Range("F133").Value = "Puntarenas"
$ws.Range("G133").Value = "Herediano"
$ws.Range("H133").Value = 1
$ws.Range("I133").Value = 2
$ws.Range("J133").Value = "A"
$ws.Range("K133").Value = 3.5
$ws.Range("L133").Value = 3.3
$ws.Range("M133").Value = 1.909
$ws.Range("N133").Value = 4.5
$ws.Range("O133").Value = 3.6
$ws.Range("P133").Value = 1.65
$ws.Range("Q133").Value = 0.75
$ws.Range("R133").Value = 1.925
$ws.Range("S133").Value = 1.875
$ws.Range("T133").Value = 2.5
$ws.Range("U133").Value = 1.9
$ws.Range("V133").Value = 1.9
$ws.Range("W133").Value = -1
$ws.Range("X133").Value = -1
$ws.Range("Y133").Value = 0.6499999999999999
$ws.Range("Z133").Value = -0.5
$ws.Range("AA133").Value = 0.4375
$ws.Range("AB133").Value = 0.8999999999999999
$ws.Range("AC133").Value = -1

# Rows 178 <-> 179
$ws.Range("B178").Value = 6782522
$ws.Range("F178").Value = "Municipal Perez Zeledon"
$ws.Range("G178").Value = "Sporting San Jose"
$ws.Range("H178").Value = 1
$ws.Range("I178").Value = 2
$ws.Range("J178").Value = "A"
$ws.Range("K178").Value = 2.5
$ws.Range("L178").Value = 3.5
$ws.Range("M178").Value = 2.5
$ws.Range("N178").Value = 2.2
$ws.Range("O178").Value = 3.5
$ws.Range("P178").Value = 2.9
$ws.Range("Q178").Value = -0.25
$ws.Range("R178").Value = 1.9
$ws.Range("S178").Value = 1.9
$ws.Range("T178").Value = 2.5
$ws.Range("U178").Value = 1.9
$ws.Range("V178").Value = 1.9
$ws.Range("W178").Value = -1
$ws.Range("X178").Value = -1
$ws.Range("Y178").Value = 1.9
$ws.Range("Z178").Value = -1
$ws.Range("AA178").Value = 0.8999999999999999
$ws.Range("AB178").Value = 0.8999999999999999
$ws.Range("AC178").Value = -1

$ws.Range("B179").Value = 6781354
$ws.Range("F179").Value = "Puntarenas"
$ws.Range("G179").Value = "AD San Carlos"
$ws.Range("H179").Value = 1
$ws.Range("I179").Value = 0
$ws.Range("J179").Value = "H"
$ws.Range("K179").Value = 2.4
$ws.Range("L179").Value = 3.2
$ws.Range("M179").Value = 2.8
$ws.Range("N179").Value = 2.3
$ws.Range("O179").Value = 3.2
$ws.Range("P179").Value = 3
$ws.Range("Q179").Value = -0.25
$ws.Range("R179").Value = 2
$ws.Range("S179").Value = 1.8
$ws.Range("T179").Value = 2.25
$ws.Range("U179").Value = 1.9
$ws.Range("V179").Value = 1.9
$ws.Range("W179").Value = 1.3
$ws.Range("X179").Value = -1
$ws.Range("Y179").Value = -1
$ws.Range("Z179").Value = 1
$ws.Range("AA179").Value = -1
$ws.Range("AB179").Value = -1
$ws.Range("AC179").Value = 0.8999999999999999

# --- Part 2: refresh rows 318 and 319 with final match data (scores now known) ---

# Row 318
$ws.Range("A318").Value = 316
$ws.Range("B318").Value = 7623984
$ws.Range("C318").Value = "Costa Rica Primera Division"
$ws.Range("D318").Value = "Costa Rica Primera Division"
$ws.Range("E318").Value = 45338.875
$ws.Range("F318").Value = "Sporting San Jose"
$ws.Range("G318").Value = "AD San Carlos"
$ws.Range("H318").Value = 0
$ws.Range("I318").Value = 2
$ws.Range("J318").Value = "A"
$ws.Range("K318").Value = 2.7
$ws.Range("L318").Value = 3.25
$ws.Range("M318").Value = 2.3
$ws.Range("N318").Value = 2.4
$ws.Range("O318").Value = 3.25
$ws.Range("P318").Value = 2.6
$ws.Range("Q318").Value = 0
$ws.Range("R318").Value = 1.85
$ws.Range("S318").Value = 1.95
$ws.Range("T318").Value = 2.5
$ws.Range("U318").Value = 1.85
$ws.Range("V318").Value = 1.95
$ws.Range("W318").Value = -1
$ws.Range("X318").Value = -1
$ws.Range("Y318").Value = 1.6
$ws.Range("Z318").Value = -1
$ws.Range("AA318").Value = 0.95
$ws.Range("AB318").Value = -1
$ws.Range("AC318").Value = 0.95

# Row 319
$ws.Range("A319").Value = 317
$ws.Range("B319").Value = 7623913
$ws.Range("C319").Value = "Costa Rica Primera Division"
$ws.Range("D319").Value = "Costa Rica Primera Division"
$ws.Range("E319").Value = 45339
$ws.Range("F319").Value = "Deportivo Saprissa"
$ws.Range("G319").Value = "Alajuelense"
$ws.Range("H319").Value = 0
$ws.Range("I319").Value = 0
$ws.Range("J319").Value = "D"
$ws.Range("K319").Value = 2.25
$ws.Range("L319").Value = 3.1
$ws.Range("M319").Value = 2.875
$ws.Range("N319").Value = 2.1
$ws.Range("O319").Value = 3.1
$ws.Range("P319").Value = 3.2
$ws.Range("Q319").Value = -0.25
$ws.Range("R319").Value = 1.9
$ws.Range("S319").Value = 1.9
$ws.Range("T319").Value = 2.25
$ws.Range("U319").Value = 1.85
$ws.Range("V319").Value = 1.95
$ws.Range("W319").Value = -1
$ws.Range("X319").Value = 2.1
$ws.Range("Y319").Value = -1
$ws.Range("Z319").Value = -0.5
$ws.Range("AA319").Value = 0.45
$ws.Range("AB319").Value = -1
$ws.Range("AC319").Value = 0.95

# --- Part 3: append new rows 320-322 (new matches), copying number format from row 319 ---

$ws.Range("A319:AC319").Copy() | Out-Null
$ws.Range("A320:AC320").PasteSpecial(-4122) | Out-Null
$ws.Range("A320").Value = 318
$ws.Range("B320").Value = 7623914
$ws.Range("C320").Value = "Costa Rica Primera Division"
$ws.Range("D320").Value = "Costa Rica Primera Division"
$ws.Range("E320").Value = 45339.75
$ws.Range("F320").Value = "AD Grecia"
$ws.Range("G320").Value = "Municipal Perez Zeledon"
$ws.Range("H320").Value = 2
$ws.Range("I320").Value = 1
$ws.Range("J320").Value = "H"
$ws.Range("K320").Value = 2
$ws.Range("L320").Value = 3.25
$ws.Range("M320").Value = 3.3
$ws.Range("N320").Value = 2.25
$ws.Range("O320").Value = 3.1
$ws.Range("P320").Value = 2.9
$ws.Range("Q320").Value = -0.25
$ws.Range("R320").Value = 2
$ws.Range("S320").Value = 1.8
$ws.Range("T320").Value = 2.25
$ws.Range("U320").Value = 1.8
$ws.Range("V320").Value = 2
$ws.Range("W320").Value = 1.25
$ws.Range("X320").Value = -1
$ws.Range("Y320").Value = -1
$ws.Range("Z320").Value = 1
$ws.Range("AA320").Value = -1
$ws.Range("AB320").Value = 0.8
$ws.Range("AC320").Value = -1

$ws.Range("A319:AC319").Copy() | Out-Null
$ws.Range("A321:AC321").PasteSpecial(-4122) | Out-Null
$ws.Range("A321").Value = 319
$ws.Range("B321").Value = 7623915
$ws.Range("C321").Value = "Costa Rica Primera Division"
$ws.Range("D321").Value = "Costa Rica Primera Division"
$ws.Range("E321").Value = 45339.92708333334
$ws.Range("F321").Value = "Puntarenas"
$ws.Range("G321").Value = "Santos de Gupiles"
$ws.Range("H321").Value = 3
$ws.Range("I321").Value = 1
$ws.Range("J321").Value = "H"
$ws.Range("K321").Value = 2
$ws.Range("L321").Value = 3.4
$ws.Range("M321").Value = 3.2
$ws.Range("N321").Value = 1.727
$ws.Range("O321").Value = 3.5
$ws.Range("P321").Value = 4.2
$ws.Range("Q321").Value = -0.75
$ws.Range("R321").Value = 2
$ws.Range("S321").Value = 1.8
$ws.Range("T321").Value = 2.5
$ws.Range("U321").Value = 1.95
$ws.Range("V321").Value = 1.85
$ws.Range("W321").Value = 0.7270000000000001
$ws.Range("X321").Value = -1
$ws.Range("Y321").Value = -1
$ws.Range("Z321").Value = 1
$ws.Range("AA321").Value = -1
$ws.Range("AB321").Value = 0.95
$ws.Range("AC321").Value = -1

$ws.Range("A319:AC319").Copy() | Out-Null
$ws.Range("A322:AC322").PasteSpecial(-4122) | Out-Null
$ws.Range("A322").Value = 320
$ws.Range("B322").Value = 7623911
$ws.Range("C322").Value = "Costa Rica Primera Division"
$ws.Range("D322").Value = "Costa Rica Primera Division"
$ws.Range("E322").Value = 45340.58333333334
$ws.Range("F322").Value = "Cartagines"
$ws.Range("G322").Value = "Municipal Liberia"
$ws.Range("H322").Value = 0
$ws.Range("I322").Value = 1
$ws.Range("J322").Value = "A"
$ws.Range("K322").Value = 1.666
$ws.Range("L322").Value = 3.75
$ws.Range("M322").Value = 4.2
$ws.Range("N322").Value = 1.533
$ws.Range("O322").Value = 4
$ws.Range("P322").Value = 5.25
$ws.Range("Q322").Value = -0.75
$ws.Range("R322").Value = 1.7
$ws.Range("S322").Value = 2.1
$ws.Range("T322").Value = 2.75
$ws.Range("U322").Value = 1.825
$ws.Range("V322").Value = 1.975
$ws.Range("W322").Value = -1
$ws.Range("X322").Value = -1
$ws.Range("Y322").Value = 4.25
$ws.Range("Z322").Value = -1
$ws.Range("AA322").Value = 1.1
$ws.Range("AB322").Value = -1
$ws.Range("AC322").Value = 0.9750000000000001

$excel.CutCopyMode = 0

# --- Part 4: update the used-range dimension to reflect the 3 new rows ---
Write-Host "Edit applied."